# Apply scheduled-runner updates to Sargatanas_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 35.347828
$ws.Range("H33").Value = 1474.9166
$ws.Range("H46").Value = 5000
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 15000
$ws.Range("N46").Value = -15238
$ws.Range("H60").Value = 5000
$ws.Range("J60").Value = 5000
$ws.Range("L60").Value = 15000
$ws.Range("N60").Value = -15968
$ws.Range("H99").Value = 339
$ws.Range("I99").Value = 339
$ws.Range("K99").Value = 1017
$ws.Range("M99").Value = 481
$ws.Range("H137").Value = 2467.28
$ws.Range("I137").Value = 2441.8635
$ws.Range("K137").Value = 7325.5905
$ws.Range("M137").Value = -4775.5905
$ws.Range("H138").Value = 4922.476
$ws.Range("I138").Value = 1985
$ws.Range("J138").Value = 5349.7456
$ws.Range("K138").Value = 5955
$ws.Range("L138").Value = 16049.2368
$ws.Range("M138").Value = -815
$ws.Range("N138").Value = -26329.2368
$ws.Range("H141").Value = 3998.5
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3231.14
$ws.Range("I32").Value = 3322.0938
$ws.Range("K32").Value = 3322.0938
$ws.Range("M32").Value = -3035.0938
$ws.Range("H122").Value = 11741.038
$ws.Range("I122").Value = 14662.6875
$ws.Range("J122").Value = 7066.4
$ws.Range("K122").Value = 43988.0625
$ws.Range("L122").Value = 21199.2
$ws.Range("M122").Value = -41538.0625
$ws.Range("N122").Value = -26099.2
$ws.Range("H128").Value = 60000
$ws.Range("J128").Value = 60000
$ws.Range("L128").Value = 60000
$ws.Range("N128").Value = -69960

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 35714680
$ws.Range("J80").Value = 341.875
$ws.Range("L80").Value = 341.875
$ws.Range("N80").Value = -2337.875
$ws.Range("H83").Value = 35714680
$ws.Range("J83").Value = 341.875
$ws.Range("L83").Value = 1709.375
$ws.Range("N83").Value = -11693.375
$ws.Range("H105").Value = 3844.8462
$ws.Range("I105").Value = 2983.7144
$ws.Range("J105").Value = 4849.5
$ws.Range("K105").Value = 2983.7144
$ws.Range("L105").Value = 4849.5
$ws.Range("M105").Value = -1236.7144
$ws.Range("N105").Value = -8343.5
$ws.Range("H107").Value = 86542424
$ws.Range("I107").Value = 140625890
$ws.Range("K107").Value = 140625890
$ws.Range("M107").Value = -140623970
$ws.Range("H134").Value = 6684.9
$ws.Range("I134").Value = 2516.1538
$ws.Range("K134").Value = 7548.4614
$ws.Range("M134").Value = -5013.4614
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5558.8594
$ws.Range("I31").Value = 2222.4783
$ws.Range("K31").Value = 2222.4783
$ws.Range("M31").Value = -1927.4783
$ws.Range("H34").Value = 5558.8594
$ws.Range("I34").Value = 2222.4783
$ws.Range("K34").Value = 2222.4783
$ws.Range("M34").Value = -2020.4783
$ws.Range("H58").Value = 8201102.5
$ws.Range("I58").Value = 11629252
$ws.Range("J58").Value = 11635.389
$ws.Range("K58").Value = 11629252
$ws.Range("L58").Value = 11635.389
$ws.Range("M58").Value = -11629049
$ws.Range("N58").Value = -12041.389
$ws.Range("H134").Value = 5061.0176
$ws.Range("I134").Value = 2226.2727
$ws.Range("K134").Value = 6678.8181
$ws.Range("M134").Value = -4143.8181
$ws.Range("H136").Value = 8201102.5
$ws.Range("I136").Value = 11629252
$ws.Range("J136").Value = 11635.389
$ws.Range("K136").Value = 34887756
$ws.Range("L136").Value = 34906.167
$ws.Range("M136").Value = -34885206
$ws.Range("N136").Value = -40006.167

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 400061920
$ws.Range("I4").Value = 400061920
$ws.Range("K4").Value = 1200185760
$ws.Range("M4").Value = -1200185648
$ws.Range("H5").Value = 1741530.4
$ws.Range("I5").Value = 3333883.5
$ws.Range("K5").Value = 10001650.5
$ws.Range("M5").Value = -10001538.5
$ws.Range("H23").Value = 274.9375
$ws.Range("J23").Value = 339.77777
$ws.Range("L23").Value = 1019.33331
$ws.Range("N23").Value = -1489.33331
$ws.Range("H86").Value = 878.2222
$ws.Range("I86").Value = 819.8
$ws.Range("K86").Value = 2459.4
$ws.Range("M86").Value = -1273.4
$ws.Range("H89").Value = 878.2222
$ws.Range("I89").Value = 819.8
$ws.Range("K89").Value = 7378.2
$ws.Range("M89").Value = -1450.2
$ws.Range("H113").Value = 1408.9
$ws.Range("I113").Value = 1220.125
$ws.Range("J113").Value = 1534.75
$ws.Range("K113").Value = 3660.375
$ws.Range("L113").Value = 4604.25
$ws.Range("M113").Value = -1490.375
$ws.Range("N113").Value = -8944.25
$ws.Range("H131").Value = 1641
$ws.Range("I131").Value = 1512
$ws.Range("K131").Value = 4536
$ws.Range("M131").Value = 504
$ws.Range("H135").Value = 1741530.4
$ws.Range("I135").Value = 3333883.5
$ws.Range("K135").Value = 30004951.5
$ws.Range("M135").Value = -30002416.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 2400
$ws.Range("I14").Value = 2400
$ws.Range("K14").Value = 2400
$ws.Range("M14").Value = -2232

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6911.8
$ws.Range("I7").Value = 4942
$ws.Range("J7").Value = 8523.454
$ws.Range("K7").Value = 4942
$ws.Range("L7").Value = 8523.454
$ws.Range("M7").Value = -4830
$ws.Range("N7").Value = -8747.454
$ws.Range("H16").Value = 1048.2
$ws.Range("I16").Value = 1122.75
$ws.Range("J16").Value = 750
$ws.Range("K16").Value = 1122.75
$ws.Range("L16").Value = 750
$ws.Range("M16").Value = -952.75
$ws.Range("N16").Value = -1090
$ws.Range("H22").Value = 3082
$ws.Range("H27").Value = 3082
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H46").Value = 7939163
$ws.Range("I46").Value = 787.25
$ws.Range("J46").Value = 11114513
$ws.Range("K46").Value = 787.25
$ws.Range("L46").Value = 11114513
$ws.Range("M46").Value = -599.25
$ws.Range("N46").Value = -11114889
$ws.Range("H93").Value = 7803.154
$ws.Range("J93").Value = 10649.833
$ws.Range("L93").Value = 10649.833
$ws.Range("N93").Value = -13145.833
$ws.Range("H122").Value = 4923.6343
$ws.Range("I122").Value = 3609.125
$ws.Range("K122").Value = 10827.375
$ws.Range("M122").Value = -8377.375
$ws.Range("H126").Value = 6911.8
$ws.Range("I126").Value = 4942
$ws.Range("J126").Value = 8523.454
$ws.Range("K126").Value = 14826
$ws.Range("L126").Value = 25570.362
$ws.Range("M126").Value = -12356
$ws.Range("N126").Value = -30510.362
$ws.Range("H132").Value = 9810199
$ws.Range("J132").Value = 9601.308000000001
$ws.Range("L132").Value = 28803.924
$ws.Range("N132").Value = -33863.924

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 60129.89
$ws.Range("J62").Value = 4635
$ws.Range("L62").Value = 4635
$ws.Range("N62").Value = -5883
$ws.Range("H65").Value = 60129.89
$ws.Range("J65").Value = 4635
$ws.Range("L65").Value = 23175
$ws.Range("N65").Value = -29415
$ws.Range("H69").Value = 22246.334
$ws.Range("I69").Value = 6749
$ws.Range("J69").Value = 29995
$ws.Range("K69").Value = 6749
$ws.Range("L69").Value = 29995
$ws.Range("M69").Value = -6000
$ws.Range("N69").Value = -31493
$ws.Range("H72").Value = 22246.334
$ws.Range("I72").Value = 6749
$ws.Range("J72").Value = 29995
$ws.Range("K72").Value = 20247
$ws.Range("L72").Value = 89985
$ws.Range("M72").Value = -16503
$ws.Range("N72").Value = -97473
$ws.Range("H126").Value = 4277
$ws.Range("I126").Value = 2582.3333
$ws.Range("J126").Value = 7666.3335
$ws.Range("K126").Value = 7746.999899999999
$ws.Range("L126").Value = 22999.0005
$ws.Range("M126").Value = -5276.999899999999
$ws.Range("N126").Value = -27939.0005
$ws.Range("H132").Value = 12205693
$ws.Range("I132").Value = 14710146
$ws.Range("K132").Value = 44130438
$ws.Range("M132").Value = -44127908
$ws.Range("H138").Value = 110000
$ws.Range("J138").Value = 110000
$ws.Range("L138").Value = 110000
$ws.Range("N138").Value = -120280

